$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (before) values for the rows that will be rotated.
# New row 2 <- old row 4, new row 3 <- old row 2, new row 4 <- old row 3
# (i.e. the three data rows end up sorted by date descending)
$rows = @(2, 3, 4)
$cols = @("D", "J", "K", "L", "M", "O", "P")

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

$mapping = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $orig[$srcRow][$c]
    }
}
